$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# A handful of Price cells (column D) hold plain numeric-looking text
# (e.g. "579.37") that Excel would otherwise auto-convert to a Number
# on assignment; a leading apostrophe keeps them literal text, matching
# how the original cell was authored, without touching cell formatting.
$ws.Range("D2").Value = "64.947.94"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "3.144.99"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D5").Value = "'579.37"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "'148.53"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.143.95"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'0.498"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "'0.0000263"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'37.02"
$ws.Range("E14").Value = "  -3.38%  "
$ws.Range("D15").Value = "3.664.09"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "64.907.08"
$ws.Range("D17").Value = "3.152.15"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "'7.12"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'502.27"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").Value = "'15.03"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "'0.712"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("D24").Value = "'7.71"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").Value = "'84.08"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'9.05"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").Value = "'27.45"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "'6.38"
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("D35").Value = "'6.47"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'54.87"
$ws.Range("D37").Value = "'0.0883"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D38").Value = "'474.22"
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("D39").Value = "'0.0412"
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("D41").Value = "'8.72"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").Value = "2.995.18"
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("D43").Value = "'0.117"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").Value = "'0.281"
$ws.Range("E44").Value = "  -4.27%  "
$ws.Range("D45").Value = "'2.41"
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("D46").Value = "'28.20"
$ws.Range("E46").Value = "  -5.41%  "
$ws.Range("D47").Value = "0.0₃0590"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").Value = "'2.24"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "'2.49"
$ws.Range("E51").Value = "  +14.17%  "
